# 68_scenecat_memory_bedrooms_2.xlsx
#
# Commit: "new input files generation for english version; make only 20
# different versions and duplicate many times for 1000 subjects"
#
# The stimulus rows (2-42, one row per trial) get re-shuffled into a new
# trial order. Columns A-G (subject_id/task/block_total/block_scene/
# trial_block/trial_total/target_cat) are the fixed trial-sequence
# scaffolding and stay put; columns H-V (category/cond_cat/cond_mem/
# correct_answer/stimulus/conceptual/perceptual/typicality/n/p_*/r_*)
# travel together with the stimulus they describe, landing on a new row.
#
# The single catch-trial row also moves: it used to live at row 29
# (stimulus stimuli/catch_14.jpg) and ends up at row 15 with a new
# stimulus file (stimuli/catch_05.jpg) and blanked-out score columns;
# row 29's old slot is filled by the real trial that used to sit at row 39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that belong to / travel with a stimulus row.
$cols = @("H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

# destRow -> sourceRow giving the row whose H..V content should end up at
# destRow. Row 15 (-> $null) becomes the new catch row instead of copying
# another row's content.
$rowMap = [ordered]@{
    2  = 40; 3  = 16; 4  = 8;  5  = 13; 6  = 10
    7  = 3;  8  = 30; 9  = 34; 10 = 27; 11 = 36
    12 = 42; 13 = 33; 14 = 28; 15 = $null; 16 = 21
    17 = 5;  18 = 18; 19 = 23; 20 = 12; 21 = 7
    22 = 25; 23 = 24; 24 = 38; 25 = 11; 26 = 35
    27 = 15; 28 = 31; 29 = 39; 30 = 14; 31 = 19
    32 = 32; 33 = 4;  34 = 6;  35 = 37; 36 = 2
    37 = 9;  38 = 41; 39 = 26; 40 = 20; 41 = 17
    42 = 22
}

# Snapshot every source row's H..V values BEFORE any writes happen, since
# several destinations read from rows that are themselves overwritten
# later in the pass (e.g. row 2 <- row 40, but row 40 <- row 20, etc.).
$snapshot = @{}
foreach ($r in 2..42) {
    $row = @{}
    foreach ($c in $cols) {
        $row[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $row
}

foreach ($dest in $rowMap.Keys) {
    $src = $rowMap[$dest]
    if ($null -eq $src) { continue }
    $srcRow = $snapshot[$src]
    foreach ($c in $cols) {
        $val = $srcRow[$c]
        if ($null -eq $val) {
            $ws.Range("$c$dest").ClearContents()
        } else {
            $ws.Range("$c$dest").Value = $val
        }
    }
}

# Row 15 becomes the new catch trial: no category/cond_cat, cond_mem =
# "catch", correct_answer stays "f", a new catch stimulus file, and all
# score columns (M-V) blanked out.
$ws.Range("H15").ClearContents()
$ws.Range("I15").ClearContents()
$ws.Range("J15").Value = "catch"
$ws.Range("K15").Value = "f"
$ws.Range("L15").Value = "stimuli/catch_05.jpg"
foreach ($c in @("M","N","O","P","Q","R","S","T","U","V")) {
    $ws.Range("$c" + "15").ClearContents()
}

Write-Output "done"
